$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the Date value (row 8, column B)
$ws.Cells.Item(8, 2).Value2 = "2024-10-02T15:04:17+00:00"

# Update the Contact value (row 10, column B)
$ws.Cells.Item(10, 2).Value2 = "Ferlab.bio (http://example.org/example-publisher)"

# Insert a new row after "Contact" (row 10) for "Jurisdiction" with an empty value
$ws.Rows.Item(11).Insert()

# Match the formatting used by the rest of the data rows (copy format from row 12)
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(11, 1).Value2 = "Jurisdiction"
$ws.Cells.Item(11, 2).Value2 = ""
